$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3127251.5
$ws.Range("I113").Value = 5265213
$ws.Range("J113").Value = 2538.4614
$ws.Range("K113").Value = 5265213
$ws.Range("L113").Value = 2538.4614
$ws.Range("M113").Value = -5261959
$ws.Range("N113").Value = -9046.4614
$ws.Range("I132").Value = 3922.3
$ws.Range("K132").Value = 11766.9
$ws.Range("M132").Value = -9236.900000000001
$ws.Range("H135").Value = 38464076
$ws.Range("I135").Value = 1488.6666
$ws.Range("J135").Value = 71432010
$ws.Range("K135").Value = 13397.9994
$ws.Range("L135").Value = 642888090
$ws.Range("M135").Value = -10862.9994
$ws.Range("N135").Value = -642893160
$ws.Range("H140").Value = 68893.336
$ws.Range("J140").Value = 68893.336
$ws.Range("L140").Value = 68893.336
$ws.Range("N140").Value = -79253.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 125
$ws.Range("I4").Value = 125
$ws.Range("K4").Value = 125
$ws.Range("M4").Value = -9
$ws.Range("H37").Value = 11711.929
$ws.Range("J37").Value = 15629.125
$ws.Range("L37").Value = 15629.125
$ws.Range("N37").Value = -16175.125
$ws.Range("H80").Value = 20257.072
$ws.Range("J80").Value = 20257.072
$ws.Range("L80").Value = 20257.072
$ws.Range("N80").Value = -22253.072
$ws.Range("H83").Value = 20257.072
$ws.Range("J83").Value = 20257.072
$ws.Range("L83").Value = 60771.216
$ws.Range("N83").Value = -70755.216

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24573.818
$ws.Range("I82").Value = 3980
$ws.Range("J82").Value = 26633.2
$ws.Range("K82").Value = 3980
$ws.Range("L82").Value = 26633.2
$ws.Range("M82").Value = -3597
$ws.Range("N82").Value = -27399.2
$ws.Range("H85").Value = 24573.818
$ws.Range("I85").Value = 3980
$ws.Range("J85").Value = 26633.2
$ws.Range("K85").Value = 3980
$ws.Range("L85").Value = 26633.2
$ws.Range("M85").Value = -2654
$ws.Range("N85").Value = -29285.2
$ws.Range("H134").Value = 19151752
$ws.Range("I134").Value = 20834616
$ws.Range("J134").Value = 9054571
$ws.Range("K134").Value = 62503848
$ws.Range("L134").Value = 27163713
$ws.Range("M134").Value = -62501313
$ws.Range("N134").Value = -27168783

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 737.3333
$ws.Range("I10").Value = 106
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 106
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 33
$ws.Range("N10").Value = -2278
$ws.Range("H50").Value = 12764.333
$ws.Range("J50").Value = 12764.333
$ws.Range("L50").Value = 12764.333
$ws.Range("N50").Value = -14014.333
$ws.Range("H59").Value = 16764
$ws.Range("J59").Value = 16764
$ws.Range("L59").Value = 16764
$ws.Range("N59").Value = -19054
$ws.Range("H60").Value = 8396.4
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 8396.4
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 8396.4
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -9418.4
$ws.Range("H68").Value = 18363.363
$ws.Range("I68").Value = 26999.5
$ws.Range("J68").Value = 16444.223
$ws.Range("K68").Value = 26999.5
$ws.Range("L68").Value = 16444.223
$ws.Range("M68").Value = -26250.5
$ws.Range("N68").Value = -17942.223
$ws.Range("H71").Value = 18363.363
$ws.Range("I71").Value = 26999.5
$ws.Range("J71").Value = 16444.223
$ws.Range("K71").Value = 80998.5
$ws.Range("L71").Value = 49332.66900000001
$ws.Range("M71").Value = -77254.5
$ws.Range("N71").Value = -56820.66900000001
$ws.Range("H74").Value = 16283.833
$ws.Range("J74").Value = 16283.833
$ws.Range("L74").Value = 16283.833
$ws.Range("N74").Value = -18031.833
$ws.Range("H77").Value = 16283.833
$ws.Range("J77").Value = 16283.833
$ws.Range("L77").Value = 48851.499
$ws.Range("N77").Value = -57587.499
$ws.Range("H132").Value = 3672.8333
$ws.Range("I132").Value = 3887.3333
$ws.Range("J132").Value = 3601.3333
$ws.Range("K132").Value = 11661.9999
$ws.Range("L132").Value = 10803.9999
$ws.Range("M132").Value = -9131.999899999999
$ws.Range("N132").Value = -15863.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1035.75
$ws.Range("J4").Value = 2339.4
$ws.Range("L4").Value = 7018.200000000001
$ws.Range("N4").Value = -7242.200000000001
$ws.Range("H113").Value = 3939
$ws.Range("I113").Value = 1293.4667
$ws.Range("J113").Value = 5356.25
$ws.Range("K113").Value = 3880.4001
$ws.Range("L113").Value = 16068.75
$ws.Range("M113").Value = -1710.4001
$ws.Range("N113").Value = -20408.75
$ws.Range("H122").Value = 410.26086
$ws.Range("I122").Value = 277.05554
$ws.Range("J122").Value = 889.8
$ws.Range("K122").Value = 2493.49986
$ws.Range("L122").Value = 8008.2
$ws.Range("M122").Value = -43.4998599999999
$ws.Range("N122").Value = -12908.2
$ws.Range("H132").Value = 1735
$ws.Range("I132").Value = 1080.375
$ws.Range("J132").Value = 2084.1333
$ws.Range("K132").Value = 9723.375
$ws.Range("L132").Value = 18757.1997
$ws.Range("M132").Value = -7193.375
$ws.Range("N132").Value = -23817.1997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 8850.272000000001
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 9335.299999999999
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 9335.299999999999
$ws.Range("M46").Value = -3844
$ws.Range("N46").Value = -9647.299999999999
$ws.Range("H80").Value = 13302.777
$ws.Range("I80").Value = 6392.857
$ws.Range("J80").Value = 17700
$ws.Range("K80").Value = 6392.857
$ws.Range("L80").Value = 17700
$ws.Range("M80").Value = -5394.857
$ws.Range("N80").Value = -19696
$ws.Range("H83").Value = 13302.777
$ws.Range("I83").Value = 6392.857
$ws.Range("J83").Value = 17700
$ws.Range("K83").Value = 31964.285
$ws.Range("L83").Value = 88500
$ws.Range("M83").Value = -26972.285
$ws.Range("N83").Value = -98484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2667.3333
$ws.Range("I68").Value = 2237.6
$ws.Range("J68").Value = 2974.2856
$ws.Range("K68").Value = 2237.6
$ws.Range("L68").Value = 2974.2856
$ws.Range("M68").Value = -1488.6
$ws.Range("N68").Value = -4472.2856
$ws.Range("H71").Value = 2667.3333
$ws.Range("I71").Value = 2237.6
$ws.Range("J71").Value = 2974.2856
$ws.Range("K71").Value = 11188
$ws.Range("L71").Value = 14871.428
$ws.Range("M71").Value = -7444
$ws.Range("N71").Value = -22359.428
$ws.Range("H140").Value = 59612.43
$ws.Range("J140").Value = 56714.5
$ws.Range("L140").Value = 56714.5
$ws.Range("N140").Value = -67074.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 13625
$ws.Range("I9").Value = 13625
$ws.Range("K9").Value = 13625
$ws.Range("M9").Value = -13485
